$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Delete the "其他有價證券" (other securities) worksheet entirely.
# ---------------------------------------------------------------------
$other = $wb.Worksheets.Item("其他有價證券")
$other.Delete()

# ---------------------------------------------------------------------
# 2) Rebuild the "基金受益憑證" (fund) worksheet so it follows the same
#    column schema as the other property sheets: a proper header row
#    (name/owner/dealer/quantity/face_value/currency/total/
#    property_category/category/date/legislator_name/legislator_id/
#    source_file/index) plus the matching data row with the extra
#    metadata columns appended.
# ---------------------------------------------------------------------
$fund = $wb.Worksheets.Item("基金受益憑證")

# -- Header row (row 1) --------------------------------------------------
$fund.Range("B1").Value = "name"
$fund.Range("C1").Value = "owner"
$fund.Range("D1").Value = "dealer"
$fund.Range("E1").Value = "quantity"
$fund.Range("F1").Value = "face_value"
$fund.Range("G1").Value = "currency"
$fund.Range("H1").Value = "total"
$fund.Range("I1").Value = "property_category"
$fund.Range("J1").Value = "category"
$fund.Range("K1").Value = "date"
$fund.Range("L1").Value = "legislator_name"
$fund.Range("M1").Value = "legislator_id"
$fund.Range("N1").Value = "source_file"
$fund.Range("O1").Value = "index"

# Apply the bold/bordered header format (copied from an existing header
# cell) to every header cell so no new style entries are introduced.
$headerFormat = $fund.Range("B1")
$headerFormat.Copy()
$fund.Range("B1:O1").PasteSpecial(-4122)

# -- Data row (row 2) -----------------------------------------------------
$fund.Range("B2").Value = "第一金店頭市場基金"
$fund.Range("C2").Value = "饒月琴"
$fund.Range("D2").Value = "第一金證券投資信託股份有限公司"
$fund.Range("E2").Value = 10000
$fund.Range("F2").Value = 10
$fund.Range("G2").Value = "新臺幣"
$fund.Range("H2").Value = 100000
$fund.Range("I2").Value = "fund"
$fund.Range("J2").Value = "normal"
$fund.Range("K2").Value = "2012-04-23"
$fund.Range("L2").Value = "許忠信"
$fund.Range("M2").Value = 1749
$fund.Range("N2").Value = "tmp50641"
$fund.Range("O2").Value = 96

# Apply the plain data format (copied from an existing data cell) to the
# newly added data cells.
$dataFormat = $fund.Range("B2")
$dataFormat.Copy()
$fund.Range("I2:O2").PasteSpecial(-4122)

$excel.CutCopyMode = $false
